$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 113, shifting existing rows 113-124 down to 115-126.
$ws.Range("A113:A114").EntireRow.Insert()

# Fixed column values shared by every record in this block.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"

# New row 113: Primera quality entry for date 44449 ($/caja 12 unidades)
$ws.Cells.Item(113, 1).Value = $mercadoId
$ws.Cells.Item(113, 2).Value = $mercado
$ws.Cells.Item(113, 3).Value = $region
$ws.Cells.Item(113, 4).Value = 44449
$ws.Cells.Item(113, 5).Value = $codreg
$ws.Cells.Item(113, 6).Value = $tipo
$ws.Cells.Item(113, 7).Value = $productoId
$ws.Cells.Item(113, 8).Value = $producto
$ws.Cells.Item(113, 9).Value = $categoriaId
$ws.Cells.Item(113, 10).Value = $categoria
$ws.Cells.Item(113, 11).Value = $variedad
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 60
$ws.Cells.Item(113, 14).Value = 19000
$ws.Cells.Item(113, 15).Value = 20000
$ws.Cells.Item(113, 16).Value = 19500
$ws.Cells.Item(113, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(113, 18).Value = $origen
$ws.Cells.Item(113, 19).Value = 1625
$ws.Cells.Item(113, 20).Value = 12

# New row 114: Segunda quality entry for date 44449 ($/caja 14 unidades)
$ws.Cells.Item(114, 1).Value = $mercadoId
$ws.Cells.Item(114, 2).Value = $mercado
$ws.Cells.Item(114, 3).Value = $region
$ws.Cells.Item(114, 4).Value = 44449
$ws.Cells.Item(114, 5).Value = $codreg
$ws.Cells.Item(114, 6).Value = $tipo
$ws.Cells.Item(114, 7).Value = $productoId
$ws.Cells.Item(114, 8).Value = $producto
$ws.Cells.Item(114, 9).Value = $categoriaId
$ws.Cells.Item(114, 10).Value = $categoria
$ws.Cells.Item(114, 11).Value = $variedad
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 60
$ws.Cells.Item(114, 14).Value = 19000
$ws.Cells.Item(114, 15).Value = 20000
$ws.Cells.Item(114, 16).Value = 19500
$ws.Cells.Item(114, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(114, 18).Value = $origen
$ws.Cells.Item(114, 19).Value = 1393
$ws.Cells.Item(114, 20).Value = 14
